# Update the handback report timestamps for the "23a6da5e..." row
# in the zh-cn and de-de sheets (row 3 of each table).

$wb = $excel.ActiveWorkbook

$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("E3").Value = "2016-03-23 00:42:33"
$wsZhCn.Range("H3").Value = "2016-03-23 00:42:58"

$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("E3").Value = "2016-03-23 00:42:37"
$wsDeDe.Range("H3").Value = "2016-03-23 00:43:05"
